# CIERRE 20 OCT 23
# Advance the payroll receipt workbook from "SEMANA 41" (09-15 Oct 2023)
# to "SEMANA 42" (16-22 Oct 2023): update the week banner, the totals that
# changed for the new week (K4 and E25, whose dependent SUM formulas
# recalculate automatically), and restore the sheet's scroll/selection
# position left by the author's last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")

# Week banner text -- every other "SEMANA ..." cell on the sheet (H9, B28,
# H28, B46, H46) is a formula chained off this one, so they update for free.
$ws.Range("B9").Value = "SEMANA  42        DEL    16     Al    22    OCTUBRE    2023"

# Updated pay figures for the new week.
$ws.Range("K4").Value = 840
$ws.Range("E25").Value = 933

# Recalculate so the dependent SUM()/TODAY()-chain formulas pick up the
# new inputs (and the current clock) before saving.
$excel.Calculate()

# Restore the view state (scroll position + active selection) as left by
# the author after finishing this week's entries.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("E26").Select()
